$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 122.253015
$ws.Range("H2").Value = 366.759045
$ws.Range("I2").Value = 0.1988639364328829
$ws.Range("J2").Value = 0.1988639364328829
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 2889.547526679965
$ws.Range("R2").Value = 26005.92774011969
$ws.Range("S2").Value = 0.01357812052046301
$ws.Range("T2").Value = 0.01357812052046301
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 122.253015
$ws.Range("H3").Value = 366.759045
$ws.Range("I3").Value = 0.1988639364328829
$ws.Range("J3").Value = 0.1988639364328829
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 22163.05242499987
$ws.Range("R3").Value = 199467.4718249989
$ws.Range("S3").Value = 0.1041452317878135
$ws.Range("T3").Value = 0.1041452317878135
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 122.253015
$ws.Range("H4").Value = 366.759045
$ws.Range("I4").Value = 0.1988639364328829
$ws.Range("J4").Value = 0.1988639364328829
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 13582.37675739719
$ws.Range("R4").Value = 122241.3908165747
$ws.Range("S4").Value = 0.06382423090931934
$ws.Range("T4").Value = 0.06382423090931932
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 122.253015
$ws.Range("H5").Value = 366.759045
$ws.Range("I5").Value = 0.1988639364328829
$ws.Range("J5").Value = 0.1988639364328829
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 3685.07743976359
$ws.Range("R5").Value = 33165.69695787231
$ws.Range("S5").Value = 0.01731635321528702
$ws.Range("T5").Value = 0.01731635321528702
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 132.5447616666667
$ws.Range("H6").Value = 397.634285
$ws.Range("I6").Value = 0.2156050961899926
$ws.Range("J6").Value = 0.2156050961899926
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 3132.801168529889
$ws.Range("R6").Value = 28195.210516769
$ws.Range("S6").Value = 0.01472118089084384
$ws.Range("T6").Value = 0.01472118089084384
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 132.5447616666667
$ws.Range("H7").Value = 397.634285
$ws.Range("I7").Value = 0.2156050961899926
$ws.Range("J7").Value = 0.2156050961899926
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 24028.82662221007
$ws.Range("R7").Value = 216259.4395998906
$ws.Range("S7").Value = 0.1129125930025979
$ws.Range("T7").Value = 0.1129125930025979
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 132.5447616666667
$ws.Range("H8").Value = 397.634285
$ws.Range("I8").Value = 0.2156050961899926
$ws.Range("J8").Value = 0.2156050961899926
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 14725.79543478812
$ws.Range("R8").Value = 132532.1589130931
$ws.Range("S8").Value = 0.06919720936480814
$ws.Range("T8").Value = 0.06919720936480814
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 132.5447616666667
$ws.Range("H9").Value = 397.634285
$ws.Range("I9").Value = 0.2156050961899926
$ws.Range("J9").Value = 0.2156050961899926
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 3995.301964345625
$ws.Range("R9").Value = 35957.71767911062
$ws.Range("S9").Value = 0.01877411293174271
$ws.Range("T9").Value = 0.01877411293174271
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 320.0894206666666
$ws.Range("H10").Value = 960.2682619999999
$ws.Range("I10").Value = 0.5206762565675317
$ws.Range("J10").Value = 0.5206762565675317
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 7565.568782117884
$ws.Range("R10").Value = 68090.11903906096
$ws.Range("S10").Value = 0.03555096560307475
$ws.Range("T10").Value = 0.03555096560307475
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 320.0894206666666
$ws.Range("H11").Value = 960.2682619999999
$ws.Range("I11").Value = 0.5206762565675317
$ws.Range("J11").Value = 0.5206762565675317
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 58028.49615547862
$ws.Range("R11").Value = 522256.4653993075
$ws.Range("S11").Value = 0.2726786485237762
$ws.Range("T11").Value = 0.2726786485237762
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 320.0894206666666
$ws.Range("H12").Value = 960.2682619999999
$ws.Range("I12").Value = 0.5206762565675317
$ws.Range("J12").Value = 0.5206762565675317
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 35562.10951158682
$ws.Range("R12").Value = 320058.9856042814
$ws.Range("S12").Value = 0.1671080348918968
$ws.Range("T12").Value = 0.1671080348918968
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 320.0894206666666
$ws.Range("H13").Value = 960.2682619999999
$ws.Range("I13").Value = 0.5206762565675317
$ws.Range("J13").Value = 0.5206762565675317
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 9648.46799733921
$ws.Range("R13").Value = 86836.2119760529
$ws.Range("S13").Value = 0.04533860754878392
$ws.Range("T13").Value = 0.04533860754878392
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 39.86989333333333
$ws.Range("H14").Value = 119.60968
$ws.Range("I14").Value = 0.06485471080959287
$ws.Range("J14").Value = 0.06485471080959287
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 942.3567318182487
$ws.Range("R14").Value = 8481.210586364239
$ws.Range("S14").Value = 0.004428178861830152
$ws.Range("T14").Value = 0.004428178861830152
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 39.86989333333333
$ws.Range("H15").Value = 119.60968
$ws.Range("I15").Value = 0.06485471080959287
$ws.Range("J15").Value = 0.06485471080959287
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 7227.948825031591
$ws.Range("R15").Value = 65051.53942528432
$ws.Range("S15").Value = 0.03396447345331647
$ws.Range("T15").Value = 0.03396447345331646
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 39.86989333333333
$ws.Range("H16").Value = 119.60968
$ws.Range("I16").Value = 0.06485471080959287
$ws.Range("J16").Value = 0.06485471080959287
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 4429.566931585057
$ws.Range("R16").Value = 39866.10238426551
$ws.Range("S16").Value = 0.02081474455608803
$ws.Range("T16").Value = 0.02081474455608803
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 39.86989333333333
$ws.Range("H17").Value = 119.60968
$ws.Range("I17").Value = 0.06485471080959287
$ws.Range("J17").Value = 0.06485471080959287
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 1201.799762962471
$ws.Range("R17").Value = 1201.799762962471
$ws.Range("S17").Value = 0.005647313938358225
$ws.Range("T17").Value = 0.005647313938358225
